$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Text changes: "Brytpunkt 1" / "Brytpunkt 2" become more descriptive labels
$ws.Range("A3").Value = "Brytpunkt för statlig inkomstskatt"
$ws.Range("A4").Value = "Brytpunkt för värnskatt"

# Column A width change 24 -> 25 (raw OOXML width units; COM ColumnWidth
# differs from the stored XML width by 5/6, so compensate accordingly)
$ws.Columns.Item(1).ColumnWidth = 24.166666666666668

# Selection change
$ws.Range("A5").Select()
